$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump each year in column A (rows 2-10) by one: 2010-2018 -> 2011-2019
for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $cell.Value2 + 1
}

# Add a new totals row (row 11) with the same look as the data rows above it,
# then drop in AVERAGE() formulas for the two numeric columns.
$srcRow = $ws.Range("A10:C10")
$dstRow = $ws.Range("A11:C11")
$srcRow.Copy()
$dstRow.PasteSpecial(-4122)

$ws.Range("B11:C11").Formula = "=AVERAGE(B2:B10)"
